$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values (weather cluster re-run fixed AF / full size output)
$ws.Range("D2").Value = 0
$ws.Range("F2").Value = 0.02089864158829677
$ws.Range("G2").Value = 0.03749999999999999
$ws.Range("I2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("F3").Value = 0.06374085684430505
$ws.Range("G3").Value = 0.0361111111111111
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0.00544959128065395
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0.01623720437698552
$ws.Range("E4").Value = 0.05267175572519084
$ws.Range("H4").Value = 0.00860215053763441
$ws.Range("I4").Value = 0.1594798083504449
$ws.Range("B5").Value = 0.002221516978736909
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("F5").Value = 0.2408568443051225
$ws.Range("G5").Value = 0.2097222222222217
$ws.Range("H5").Value = 0.00860215053763441
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0.3835149863760199
$ws.Range("D6").Value = 0
$ws.Range("F6").Value = 0.01880877742946709
$ws.Range("G6").Value = 0.0125
$ws.Range("I6").Value = 0
$ws.Range("E7").Value = 0.03511450381679389
$ws.Range("B8").Value = 0.443351317042196
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0.0858947368421052
$ws.Range("F8").Value = 0.1003134796238243
$ws.Range("G8").Value = 0.06666666666666668
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0.2568119891008189
$ws.Range("E9").Value = 0.02977099236641221
$ws.Range("K9").Value = 0.01449275362318841
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0.1200141193081534
$ws.Range("E10").Value = 0.04351145038167939
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0.09856262833675584
$ws.Range("E12").Value = 0.05190839694656488
$ws.Range("F12").Value = 0.001567398119122257
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("F13").Value = 0.04493207941483803
$ws.Range("G13").Value = 0.002777777777777778
$ws.Range("I13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("F14").Value = 0.09926854754440942
$ws.Range("G14").Value = 0.481944444444443
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0.07288828337874667
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0.02012001411930813
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0.03282442748091603
$ws.Range("F16").Value = 0.02246603970741904
$ws.Range("G16").Value = 0.001388888888888889
$ws.Range("H16").Value = 0.3053763440860215
$ws.Range("I16").Value = 0.04038329911019845
$ws.Range("D17").Value = 0
$ws.Range("F17").Value = 0.002089864158829676
$ws.Range("G17").Value = 0.006944444444444445
$ws.Range("I17").Value = 0
$ws.Range("E18").Value = 0.005343511450381679
$ws.Range("H18").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("F19").Value = 0.08516196447230916
$ws.Range("G19").Value = 0.01111111111111111
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0.01294277929155313
$ws.Range("D20").Value = 0
$ws.Range("F20").Value = 0.01880877742946709
$ws.Range("G20").Value = 0.006944444444444445
$ws.Range("I20").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("F21").Value = 0.03709508881922677
$ws.Range("G21").Value = 0.01527777777777778
$ws.Range("I21").Value = 0
$ws.Range("C22").Value = 0.04094599364631121
$ws.Range("D23").Value = 0
$ws.Range("F23").Value = 0.1332288401253917
$ws.Range("G23").Value = 0.03749999999999999
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0.00544959128065395
$ws.Range("E24").Value = 0.009923664122137405
$ws.Range("H24").Value = 0
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 0.05683021531944914
$ws.Range("E32").Value = 0.02442748091603053
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0.3559139784946246
$ws.Range("I32").Value = 0.1033538672142371
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 0.03494528768090354
$ws.Range("E33").Value = 0.1557251908396955
$ws.Range("F33").Value = 0.004179728317659352
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0.2080766598220391
$ws.Range("F34").Value = 0.003134796238244514
$ws.Range("I34").Value = 0

# Remove the now-defunct joint-regime-area rows (36-40); dimension auto-updates to A1:K35
$ws.Range("A36:K40").EntireRow.Delete()

